$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) and a couple of the symbol strings in column E hold
# text values (not numbers) in the source data. Force the target cells to a
# text number format before writing so Excel keeps them as exact strings
# instead of converting to floating point numbers.

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextValue "D2"  "271.34"
Set-TextValue "D3"  "23.11"
Set-TextValue "D4"  "6.382"
Set-TextValue "D5"  "0.06306"
Set-TextValue "D6"  "3.663"
Set-TextValue "D7"  "6.770"
Set-TextValue "D8"  "1.402"
Set-TextValue "D9"  "0.8348"
Set-TextValue "D10" "0.1634"
Set-TextValue "D11" "0.08462"
Set-TextValue "D13" "0.03123"
Set-TextValue "D14" "0.09308"
Set-TextValue "D15" "3.920"
Set-TextValue "D16" "0.001721"
Set-TextValue "D18" "0.006246"
Set-TextValue "D19" "0.005476"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "3.732"
Set-TextValue "D23" "2.350"
Set-TextValue "D24" "0.01388"
Set-TextValue "D27" "0.0003743"
Set-TextValue "D41" "0.006890"
Set-TextValue "D42" "0.1176"
Set-TextValue "D43" "0.003457"
Set-TextValue "D44" "0.01258"
Set-TextValue "D45" "0.00006269"

Set-TextValue "D47" "0.6452"
Set-TextValue "E47" "46CoinbaseStockTokenCOINWorstin24h"

Set-TextValue "D48" "0.1092"

Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "48CryptobidCoinCBC"

Set-TextValue "D50" "0.01240"
